$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LINK")

# Row 21 - Earnings Before Interest And Taxes
$ws.Range("D21").Value = 2200
$ws.Range("E21").Value = 2900
$ws.Range("G21").Value = 1700
$ws.Range("I21").Value = 300
$ws.Range("J21").Value = "NA"

# Row 83 - Depreciation
$ws.Range("D83").Value = 100
$ws.Range("I83").Value = 200
$ws.Range("J83").Value = "NA"

# Row 89 - Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 1900
$ws.Range("E89").Value = 2100
$ws.Range("F89").Value = 1000
$ws.Range("G89").Value = -500
$ws.Range("H89").Value = -100
$ws.Range("I89").Value = -1200
$ws.Range("J89").Value = 100

# Row 91 - Capital Expenditures
$ws.Range("D91").Value = -200
$ws.Range("E91").Value = -100
$ws.Range("G91").Value = -100
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = -200
$ws.Range("J91").Value = -100

# Row 94 - Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -200
$ws.Range("E94").Value = -100
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = -200
$ws.Range("J94").Value = "NA"

# Row 100 - Total Cash Flows From Financing Activities
$ws.Range("D100").Value = "NA"
$ws.Range("E100").Value = 0
$ws.Range("J100").Value = "NA"

# Row 101 - Effect Of Exchange Rate Changes
$ws.Range("D101").Value = -100
$ws.Range("E101").Value = 0
$ws.Range("J101").Value = "NA"

# Row 102 - Change In Cash and Cash Equivalents
$ws.Range("D102").Value = 1600
$ws.Range("E102").Value = 2000
$ws.Range("F102").Value = 900
$ws.Range("G102").Value = -600
$ws.Range("H102").Value = -200
$ws.Range("I102").Value = -1400
$ws.Range("J102").Value = 300
